$p = $ppt.ActivePresentation

# 1. The table on slide 5 gets a new table style id.
$s5 = $p.Slides.Item(5)
foreach ($sh in $s5.Shapes) {
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{48F8F92C-94EB-4D6E-88F2-9879BA0215D1}")
    }
}

# 2. Theme swap: the deck's active theme (Integral / Red Violet colours,
#    used by the slide master + slides) is replaced with the plain
#    "Office" colour scheme that used to live only in the notes-master
#    theme part.
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
